# Bugfixed the naive forecaster component module
#
# The "date" column (A) previously held Excel date-serial numbers
# (formatted via a custom numFmt as YYYY-MM-DD HH:MM:SS) representing the
# Q4 date of each year. The naive forecaster's date handling was buggy,
# so we switch the date column over to plain quarter-label text
# (e.g. "1987Q4") matching the rest of the pipeline, and drop the now
# unused custom date number format / style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$quarterLabels = @(
    "1987Q4", "1988Q4", "1989Q4", "1990Q4", "1991Q4", "1992Q4", "1993Q4",
    "1994Q4", "1995Q4", "1996Q4", "1997Q4", "1998Q4", "1999Q4", "2000Q4",
    "2001Q4", "2002Q4", "2003Q4", "2004Q4", "2005Q4", "2006Q4", "2007Q4",
    "2008Q4", "2009Q4", "2010Q4", "2011Q4", "2012Q4", "2013Q4", "2014Q4",
    "2015Q4", "2016Q4", "2017Q4", "2018Q4", "2019Q4", "2020Q4", "2021Q4",
    "2022Q4", "2023Q4", "2024Q4"
)

# Use the same style as the header cells (A1/B1): bordered, bold, centered
# text - but without the custom date number format that used to live on
# the data rows.
$headerStyle = $ws.Range("A1").Style

$row = 2
foreach ($label in $quarterLabels) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $label
    $cell.Style = $headerStyle
    $row = $row + 1
}
